# Add "Wins", "Losses", "Ties" columns (AD, AE, AF) to the season-record
# table, populating the header row and all 46 data rows with the team's
# 2017 season record (93 wins, 69 losses, 0 ties).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Header row (row 1) --------------------------------------------------
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Reuse the same header formatting (bold, centered, thin border) used by
# the rest of row 1, without introducing a brand-new style definition.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# -- Data rows (rows 2-47) ------------------------------------------------
$lastRow = 47
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 93
    $ws.Cells.Item($r, 31).Value = 69
    $ws.Cells.Item($r, 32).Value = 0
}
